$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.246.20"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.990.30"
$ws.Range("E3").Value = "  +6.09%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5099"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4113"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08793"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("D13").Value = "1.981.06"
$ws.Range("E13").Value = "  +5.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.483"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.401"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001120"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06538"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.077"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").Value = "30.282.88"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.210"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "2.212.35"
$ws.Range("E26").Value = "  +6.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.375"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.135"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.039"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.797"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.316"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02481"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.388"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06493"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2177"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.904"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6558"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.219"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6127"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.187"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.655"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06874"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.16%  "
